$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New daily rows appended to the report (update through 13/05/2021).
$newRows = @(
    @{ Row = 252; Date = 44326; B = 0; C = 7; D = 64.486411791801 },
    @{ Row = 253; Date = 44327; B = 0; C = 6; D = 55.27406725011516 },
    @{ Row = 254; Date = 44328; B = 0; C = 6; D = 55.27406725011516 },
    @{ Row = 255; Date = 44329; B = 1; C = 5; D = 46.06172270842929 }
)

foreach ($r in $newRows) {
    # Clone the formatting (style) of the last existing data row (251) onto
    # column A of the new row, then fill in the values.
    $ws.Range("A251").Copy($ws.Range("A" + $r.Row))

    $ws.Range("A" + $r.Row).Value = $r.Date
    $ws.Range("B" + $r.Row).Value = $r.B
    $ws.Range("C" + $r.Row).Value = $r.C
    $ws.Range("D" + $r.Row).Value = $r.D
}
